# "adjust spreadsheet for spores/no spores"
#
# The "Stool Weights" sheet uses volatile RANDBETWEEN() formulas (and
# difference formulas that depend on them) to synthesize example mouse
# stool-weight data. The commit simply re-opened the workbook and forced a
# fresh recalculation (F9 / Calculate Now), which re-rolled every
# RANDBETWEEN() result - and, in turn, every "today - baseline" difference
# cell that reads from it - to new random values. The author also ended up
# with the selection sitting on the last cell of the table (Q47) instead of
# the stale saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stool Weights")

# Make it the active sheet (matches tabSelected="1" on this sheet).
$ws.Activate()

# Force a full recalculation of all volatile formulas (RANDBETWEEN, etc.)
# across the workbook, same as pressing F9 / Ctrl+Alt+F9 in Excel.
$excel.CalculateFull()

# Leave the selection on the bottom-right cell of the table, same as the
# saved state in the edited workbook.
$ws.Range("Q47").Select()
